$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# The two observation records currently stored in rows 21 and 22 were
# re-ordered upstream (same two rows, values swapped). Read every value
# that differs between the rows first, then write the swapped values back,
# so we don't clobber data we still need to read.

$cols = @("A","B","C","D","E","F","G","H","J","K","N","Q","R","AF","AJ","AK","AO")

$row21 = @{}
$row22 = @{}
foreach ($col in $cols) {
    $row21[$col] = $ws.Range($col + "21").Value2()
    $row22[$col] = $ws.Range($col + "22").Value2()
}

foreach ($col in $cols) {
    $ws.Range($col + "21").Value = $row22[$col]
    $ws.Range($col + "22").Value = $row21[$col]
}

# J21/K21/N21/AF21/AJ21/AK21/AO21 were non-empty (present-but-blank, or
# substrate info) before the swap and must end up empty on row 21 after it,
# since row 22's version of the record never had them populated.
$ws.Range("J21").ClearContents()
$ws.Range("K21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("AF21").ClearContents()
$ws.Range("AJ21").ClearContents()
$ws.Range("AK21").ClearContents()
$ws.Range("AO21").ClearContents()
